# Update LR-pair metrics with recomputed TPM-based values (rows 2-10,
# columns G-T) per updated NATMI scripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 21.50919633333334
$ws.Range("H2").Value = 64.52758900000001
$ws.Range("I2").Value = 0.1832723264758264
$ws.Range("J2").Value = 0.1832723264758264
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04016966666666667
$ws.Range("N2").Value = 0.120509
$ws.Range("O2").Value = 0.01016394183724367
$ws.Range("P2").Value = 0.01016394183724367
$ws.Range("Q2").Value = 0.864017246977889
$ws.Range("R2").Value = 7.776155222801001
$ws.Range("S2").Value = 0.001862769266676633
$ws.Range("T2").Value = 0.001862769266676633
$ws.Range("G3").Value = 21.50919633333334
$ws.Range("H3").Value = 64.52758900000001
$ws.Range("I3").Value = 0.1832723264758264
$ws.Range("J3").Value = 0.1832723264758264
$ws.Range("O3").Value = 0.1167952962934662
$ws.Range("P3").Value = 0.1167952962934662
$ws.Range("Q3").Value = 9.928544651217113
$ws.Range("R3").Value = 89.35690186095401
$ws.Range("S3").Value = 0.02140534567313701
$ws.Range("T3").Value = 0.02140534567313701
$ws.Range("G4").Value = 21.50919633333334
$ws.Range("H4").Value = 64.52758900000001
$ws.Range("I4").Value = 0.1832723264758264
$ws.Range("J4").Value = 0.1832723264758264
$ws.Range("M4").Value = 3.450409
$ws.Range("N4").Value = 10.351227
$ws.Range("O4").Value = 0.8730407618692901
$ws.Range("P4").Value = 0.8730407618692901
$ws.Range("Q4").Value = 74.21552461130034
$ws.Range("R4").Value = 667.939721501703
$ws.Range("S4").Value = 0.1600042115360128
$ws.Range("T4").Value = 0.1600042115360128
$ws.Range("I5").Value = 0.5927317426910698
$ws.Range("J5").Value = 0.5927317426910698
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04016966666666667
$ws.Range("N5").Value = 0.120509
$ws.Range("O5").Value = 0.01016394183724367
$ws.Range("P5").Value = 0.01016394183724367
$ws.Range("Q5").Value = 2.794368677280333
$ws.Range("R5").Value = 25.149318095523
$ws.Range("S5").Value = 0.006024490957800114
$ws.Range("T5").Value = 0.006024490957800115
$ws.Range("I6").Value = 0.5927317426910698
$ws.Range("J6").Value = 0.5927317426910698
$ws.Range("O6").Value = 0.1167952962934662
$ws.Range("P6").Value = 0.1167952962934662
$ws.Range("S6").Value = 0.06922827951014604
$ws.Range("T6").Value = 0.06922827951014604
$ws.Range("I7").Value = 0.5927317426910698
$ws.Range("J7").Value = 0.5927317426910698
$ws.Range("M7").Value = 3.450409
$ws.Range("N7").Value = 10.351227
$ws.Range("O7").Value = 0.8730407618692901
$ws.Range("P7").Value = 0.8730407618692901
$ws.Range("Q7").Value = 240.024765786941
$ws.Range("R7").Value = 2160.222892082469
$ws.Range("S7").Value = 0.5174789722231236
$ws.Range("T7").Value = 0.5174789722231236
$ws.Range("G8").Value = 26.28859766666667
$ws.Range("H8").Value = 78.86579300000001
$ws.Range("I8").Value = 0.2239959308331038
$ws.Range("J8").Value = 0.2239959308331038
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04016966666666667
$ws.Range("N8").Value = 0.120509
$ws.Range("O8").Value = 0.01016394183724367
$ws.Range("P8").Value = 0.01016394183724367
$ws.Range("Q8").Value = 1.056004205404111
$ws.Range("R8").Value = 9.504037848637001
$ws.Range("S8").Value = 0.002276681612766923
$ws.Range("T8").Value = 0.002276681612766923
$ws.Range("G9").Value = 26.28859766666667
$ws.Range("H9").Value = 78.86579300000001
$ws.Range("I9").Value = 0.2239959308331038
$ws.Range("J9").Value = 0.2239959308331038
$ws.Range("O9").Value = 0.1167952962934662
$ws.Range("P9").Value = 0.1167952962934662
$ws.Range("Q9").Value = 12.13469400281089
$ws.Range("R9").Value = 109.212246025298
$ws.Range("S9").Value = 0.02616167111018311
$ws.Range("T9").Value = 0.02616167111018311
$ws.Range("G10").Value = 26.28859766666667
$ws.Range("H10").Value = 78.86579300000001
$ws.Range("I10").Value = 0.2239959308331038
$ws.Range("J10").Value = 0.2239959308331038
$ws.Range("M10").Value = 3.450409
$ws.Range("N10").Value = 10.351227
$ws.Range("O10").Value = 0.8730407618692901
$ws.Range("P10").Value = 0.8730407618692901
$ws.Range("Q10").Value = 90.70641398644568
$ws.Range("R10").Value = 816.3577258780111
$ws.Range("S10").Value = 0.1955575781101538
$ws.Range("T10").Value = 0.1955575781101537